# Update the "Förändrad" (Changed) date column (C) for rows 2 through 108
# from 45224 (2023-10-25) to 45233 (2023-11-03).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C108").Value = 45233
